# Update Excel files after daily scrape - 2025-10-16 03:14:00 UTC
# Refresh the opportunity listing data (rows 2-20), adjust the PREMIUM
# highlight, and resize a few columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').NumberFormat = '@'
$ws.Range('A2').Value = '1328566'
$ws.Range('A2').Style = 'Normal'
$ws.Range('B2').Value = 'https://aiesec.org/opportunity/global-talent/1328566'
$ws.Range('C2').Value = 'HR Intern'
$ws.Range('D2').Value = 'Santiago, Región Metropolitana, Chile'
$ws.Range('E2').Value = 'No'
$ws.Range('F2').Value = '3 applicants'
$ws.Range('G2').Value = '6 - 18 Months'
$ws.Range('H2').Value = 'Boehringer Ingelheim in Chile'

$ws.Range('A3').NumberFormat = '@'
$ws.Range('A3').Value = '1328548'
$ws.Range('A3').Style = 'Normal'
$ws.Range('B3').Value = 'https://aiesec.org/opportunity/global-talent/1328548'
$ws.Range('C3').Value = 'Sales'
$ws.Range('D3').Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range('E3').Value = 'No'
$ws.Range('F3').Value = '0 applicants'
$ws.Range('G3').Value = '9 - 12 Weeks'
$ws.Range('H3').Value = 'TAR - Company'

$ws.Range('A4').NumberFormat = '@'
$ws.Range('A4').Value = '1328547'
$ws.Range('A4').Style = 'Normal'
$ws.Range('B4').Value = 'https://aiesec.org/opportunity/global-talent/1328547'
$ws.Range('C4').Value = 'SEO'
$ws.Range('D4').Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range('E4').Value = 'No'
$ws.Range('F4').Value = '0 applicants'
$ws.Range('G4').Value = '9 - 12 Weeks'
$ws.Range('H4').Value = 'TAR - Company'

$ws.Range('A5').NumberFormat = '@'
$ws.Range('A5').Value = '1328545'
$ws.Range('A5').Style = 'Normal'
$ws.Range('B5').Value = 'https://aiesec.org/opportunity/global-talent/1328545'
$ws.Range('C5').Value = 'Graphic designer'
$ws.Range('D5').Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range('E5').Value = 'No'
$ws.Range('F5').Value = '0 applicants'
$ws.Range('G5').Value = '9 - 12 Weeks'
$ws.Range('H5').Value = 'TAR - Company'

$ws.Range('A6').NumberFormat = '@'
$ws.Range('A6').Value = '1328543'
$ws.Range('A6').Style = 'Normal'
$ws.Range('B6').Value = 'https://aiesec.org/opportunity/global-talent/1328543'
$ws.Range('C6').Value = 'UI/UX design'
$ws.Range('D6').Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range('E6').Value = 'No'
$ws.Range('F6').Value = '0 applicants'
$ws.Range('G6').Value = '9 - 12 Weeks'
$ws.Range('H6').Value = 'TAR - Company'

$ws.Range('A7').NumberFormat = '@'
$ws.Range('A7').Value = '1328541'
$ws.Range('A7').Style = 'Normal'
$ws.Range('B7').Value = 'https://aiesec.org/opportunity/global-talent/1328541'
$ws.Range('C7').Value = 'Marketing'
$ws.Range('D7').Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range('E7').Value = 'No'
$ws.Range('F7').Value = '0 applicants'
$ws.Range('G7').Value = '9 - 12 Weeks'
$ws.Range('H7').Value = 'TAR - Company'

$ws.Range('A8').NumberFormat = '@'
$ws.Range('A8').Value = '1328490'
$ws.Range('A8').Style = 'Normal'
$ws.Range('B8').Value = 'https://aiesec.org/opportunity/global-talent/1328490'
$ws.Range('C8').Value = 'Sales Intern'
$ws.Range('D8').Value = 'Ümraniye, Elmalıkent, 34764 Ümraniye/İstanbul, Türkiye'
$ws.Range('E8').Value = 'No'
$ws.Range('F8').Value = '6 applicants'
$ws.Range('G8').Value = '9 - 12 Weeks'
$ws.Range('H8').Value = 'ENTES ELEKTRONİK CİHAZLAR İMALAT VE TİCARET ANONİM ŞİRKETİ'

$ws.Range('A9').NumberFormat = '@'
$ws.Range('A9').Value = '1328482'
$ws.Range('A9').Style = 'Normal'
$ws.Range('B9').Value = 'https://aiesec.org/opportunity/global-talent/1328482'
$ws.Range('C9').Value = 'Web Master'
$ws.Range('D9').Value = 'İstanbul, Türkiye'
$ws.Range('E9').Value = 'No'
$ws.Range('F9').Value = '19 applicants'
$ws.Range('G9').Value = '6 - 18 Months'
$ws.Range('H9').Value = 'CCM Giyim'

$ws.Range('A10').NumberFormat = '@'
$ws.Range('A10').Value = '1328465'
$ws.Range('A10').Style = 'Normal'
$ws.Range('B10').Value = 'https://aiesec.org/opportunity/global-talent/1328465'
$ws.Range('C10').Value = 'Sales Engineer'
$ws.Range('D10').Value = 'Cairo, Cairo Governorate, Egypt'
$ws.Range('E10').Value = 'No'
$ws.Range('F10').Value = '1 applicant'
$ws.Range('G10').Value = '9 - 12 Weeks'
$ws.Range('H10').Value = 'MechNova Engineering'

$ws.Range('A11').NumberFormat = '@'
$ws.Range('A11').Value = '1326310'
$ws.Range('A11').Style = 'Normal'
$ws.Range('B11').Value = 'https://aiesec.org/opportunity/global-talent/1326310'
$ws.Range('C11').Value = 'Back - End Developer'
$ws.Range('D11').Value = 'Glyfada, Greece'
$ws.Range('E11').Value = 'No'
$ws.Range('F11').Value = '242 applicants'
$ws.Range('G11').Value = '3 - 6 Months'
$ws.Range('H11').Value = 'Validata Software'

$ws.Range('A12').NumberFormat = '@'
$ws.Range('A12').Value = '1326041'
$ws.Range('A12').Style = 'Normal'
$ws.Range('B12').Value = 'https://aiesec.org/opportunity/global-talent/1326041'
$ws.Range('C12').Value = 'ACE Program | Spanish Financial Analyst'
$ws.Range('D12').Value = 'Thane, Maharashtra, India'
$ws.Range('E12').Value = 'Yes'
$ws.Range('F12').Value = '22 applicants'
$ws.Range('G12').Value = '6 - 18 Months'
$ws.Range('H12').Value = 'Tata Consultancy Services Ltd.'

$ws.Range('A13').NumberFormat = '@'
$ws.Range('A13').Value = '1325656'
$ws.Range('A13').Style = 'Normal'
$ws.Range('B13').Value = 'https://aiesec.org/opportunity/global-talent/1325656'
$ws.Range('C13').Value = 'Design Intern'
$ws.Range('D13').Value = 'Mumbai, Maharashtra, India'
$ws.Range('E13').Value = 'No'
$ws.Range('F13').Value = '10 applicants'
$ws.Range('G13').Value = '6 - 18 Months'
$ws.Range('H13').Value = 'Rediffusion Brand Solutions Pvt Ltd'

$ws.Range('A14').NumberFormat = '@'
$ws.Range('A14').Value = '1325604'
$ws.Range('A14').Style = 'Normal'
$ws.Range('B14').Value = 'https://aiesec.org/opportunity/global-talent/1325604'
$ws.Range('C14').Value = 'International Business & Innovation Analyst'
$ws.Range('D14').Value = '4520 Santa Maria da Feira, Portugal'
$ws.Range('E14').Value = 'No'
$ws.Range('F14').Value = '143 applicants'
$ws.Range('G14').Value = '9 - 12 Weeks'
$ws.Range('H14').Value = 'M2K Consultoria'

$ws.Range('A15').NumberFormat = '@'
$ws.Range('A15').Value = '1325594'
$ws.Range('A15').Style = 'Normal'
$ws.Range('B15').Value = 'https://aiesec.org/opportunity/global-talent/1325594'
$ws.Range('C15').Value = 'Brand Partner'
$ws.Range('D15').Value = 'Mumbai, Maharashtra, India'
$ws.Range('E15').Value = 'No'
$ws.Range('F15').Value = '27 applicants'
$ws.Range('G15').Value = '6 - 18 Months'
$ws.Range('H15').Value = 'Rediffusion Brand Solutions Pvt Ltd'

$ws.Range('A16').NumberFormat = '@'
$ws.Range('A16').Value = '1325379'
$ws.Range('A16').Style = 'Normal'
$ws.Range('B16').Value = 'https://aiesec.org/opportunity/global-talent/1325379'
$ws.Range('C16').Value = 'Software Development Intern'
$ws.Range('D16').Value = 'Athens, Greece'
$ws.Range('E16').Value = 'No'
$ws.Range('F16').Value = '113 applicants'
$ws.Range('G16').Value = '9 - 12 Weeks'
$ws.Range('H16').Value = 'Eutopians'

$ws.Range('A17').NumberFormat = '@'
$ws.Range('A17').Value = '1325378'
$ws.Range('A17').Style = 'Normal'
$ws.Range('B17').Value = 'https://aiesec.org/opportunity/global-talent/1325378'
$ws.Range('C17').Value = 'Content Creation and Social Media Marketing Intern'
$ws.Range('D17').Value = 'Athens, Greece'
$ws.Range('E17').Value = 'No'
$ws.Range('F17').Value = '48 applicants'
$ws.Range('G17').Value = '9 - 12 Weeks'
$ws.Range('H17').Value = 'Eutopians'

$ws.Range('A18').NumberFormat = '@'
$ws.Range('A18').Value = '1322690'
$ws.Range('A18').Style = 'Normal'
$ws.Range('B18').Value = 'https://aiesec.org/opportunity/global-talent/1322690'
$ws.Range('C18').Value = 'Business Development Representative'
$ws.Range('D18').Value = 'Athens, Greece'
$ws.Range('E18').Value = 'No'
$ws.Range('F18').Value = '66 applicants'
$ws.Range('G18').Value = '9 - 12 Weeks'
$ws.Range('H18').Value = 'Travelr'

$ws.Range('A19').NumberFormat = '@'
$ws.Range('A19').Value = '1307292'
$ws.Range('A19').Style = 'Normal'
$ws.Range('B19').Value = 'https://aiesec.org/opportunity/global-talent/1307292'
$ws.Range('C19').Value = 'Software Engineering Intern'
$ws.Range('D19').Value = 'Hà Nội, Việt Nam'
$ws.Range('E19').Value = 'No'
$ws.Range('F19').Value = '40 applicants'
$ws.Range('G19').Value = '9 - 12 Weeks'
$ws.Range('H19').Value = 'FPT Software'

$ws.Range('A20').NumberFormat = '@'
$ws.Range('A20').Value = '1296685'
$ws.Range('A20').Style = 'Normal'
$ws.Range('B20').Value = 'https://aiesec.org/opportunity/global-talent/1296685'
$ws.Range('C20').Value = 'PHP Developer'
$ws.Range('D20').Value = 'Đống Đa, Vietnam'
$ws.Range('E20').Value = 'No'
$ws.Range('F20').Value = '372 applicants'
$ws.Range('G20').Value = '9 - 12 Weeks'
$ws.Range('H20').Value = 'Vitex Vietnam Software Joint Stock Company'


# PREMIUM (E column) highlight upkeep:
#  - E2, E5, E6 were premium ("Yes", yellow-highlighted) and are now "No" -> clear highlight
#  - E12 was not premium and is now "Yes" -> apply highlight
# Copy cell formatting instead of touching Interior.ColorIndex directly so the
# workbook's style table stays byte-for-byte equivalent to a genuine Excel edit
# (no stray per-cell styles left behind).

# Grab the still-intact "premium" (yellow) look from E2 before it gets cleared,
# and stamp it onto E12.
$ws.Range('E2').Copy()
$ws.Range('E12').PasteSpecial(-4122)

# Now reset E2, E5, E6 back to the plain (unstyled) look using E3 as the source.
$ws.Range('E3').Copy()
$ws.Range('E2').PasteSpecial(-4122)
$ws.Range('E5').PasteSpecial(-4122)
$ws.Range('E6').PasteSpecial(-4122)

# Update column widths: C, D, F, H change; A, B, E, G stay the same.
# ColumnWidth in this COM surface adds a fixed 5/6-character padding versus the
# raw OOXML <col width>, so subtract it to land on the exact target widths.
$ws.Columns.Item(3).ColumnWidth = 52.166666666666664
$ws.Columns.Item(4).ColumnWidth = 69.16666666666667
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668
$ws.Columns.Item(8).ColumnWidth = 60.166666666666664
